$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") should carry the same
# formatting (bold font, thin border, centered/top alignment) as the
# existing header row cells (e.g. H1). Copy/PasteSpecial(formats) from
# H1 reuses the existing style (s="1") instead of minting a new one.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for the new columns I (I0) and J (IF), rows 2-18.
$iValues = @(4, 9, 7, 8, 8, 8, 5, 8, 7, 8, 9, 9, 4, 8, 8, 3, 6)
$jValues = @(5, 9, 7, 8, 9, 9, 5, 8, 8, 8, 9, 9, 4, 8, 8, 3, 6)

for ($r = 0; $r -lt $iValues.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
